$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.310.29'
$ws.Range("E2").Value = '  +0.80%  '

$ws.Range("D3").Value = '3.494.39'
$ws.Range("E3").Value = '  -0.06%  '

$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '605.65'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.66%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.89'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.83%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.616'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.98%  '

$ws.Range("D8").Value = '3.487.99'
$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("E9").Value = '  -0.08%  '

$ws.Range("E10").Value = '  +2.54%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.64'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.07%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.576'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.93%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '46.70'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.85%  '

$ws.Range("E14").Value = '  +0.04%  '

$ws.Range("D15").Value = '4.054.15'
$ws.Range("E15").Value = '  -0.20%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.26'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -6.24%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '609.07'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -9.37%  '

$ws.Range("D18").Value = '3.492.80'
$ws.Range("E18").Value = '  -0.20%  '

$ws.Range("D19").Value = '69.327.98'
$ws.Range("E19").Value = '  +0.85%  '

$ws.Range("E20").Value = '  -1.87%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.14'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.00%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.96'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -10.77%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.874'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.07%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '15.70'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.28%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '95.28'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.80%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.84'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.67%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.08%  '

$ws.Range("E28").Value = '  -2.30%  '

$ws.Range("E29").Value = '  -2.47%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.04'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.61%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.38'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.62%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.06'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.08%  '

$ws.Range("E33").Value = '  -2.09%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.85'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.90%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '553.51'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.24%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.71'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.84%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.46'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.23%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '56.60'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.40%  '

$ws.Range("E39").Value = '  +0.10%  '

$ws.Range("E40").Value = '  -4.26%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0444'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.75%  '

$ws.Range("E42").Value = '  +0.82%  '

$ws.Range("D43").Value = '3.326.86'
$ws.Range("E43").Value = '  -2.39%  '

$ws.Range("E44").Value = '  -3.81%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '32.88'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.13%  '

$ws.Range("D46").Value = '0.0₃0693'
$ws.Range("E46").Value = '  -1.06%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.59'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.07%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.85'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.67%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.128'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.47%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '135.11'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.40%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.65'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +7.44%  '

